$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1, matching the style of existing header cells (A1:E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the formatting (bold, centered, bordered) from an existing header cell
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Boolean outlier flag values for rows 2-20, columns F (KNN), G (SVM), H (RF)
$values = @(
    @(0,0,0),  # row 2
    @(0,0,0),  # row 3
    @(0,0,0),  # row 4
    @(0,0,0),  # row 5
    @(0,0,0),  # row 6
    @(1,0,0),  # row 7
    @(0,0,0),  # row 8
    @(0,0,0),  # row 9
    @(0,0,0),  # row 10
    @(0,0,0),  # row 11
    @(0,0,0),  # row 12
    @(0,0,0),  # row 13
    @(0,0,0),  # row 14
    @(1,0,0),  # row 15
    @(0,0,0),  # row 16
    @(0,0,0),  # row 17
    @(1,0,1),  # row 18
    @(0,0,0),  # row 19
    @(0,0,0)   # row 20
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowVals = $values[$i]
    $ws.Cells.Item($row, 6).Value = [bool]$rowVals[0]
    $ws.Cells.Item($row, 7).Value = [bool]$rowVals[1]
    $ws.Cells.Item($row, 8).Value = [bool]$rowVals[2]
}
